$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A326").Value = "How many curves can be plotted in GEO?"
$ws.Range("B326").Value = "llama3.2:latest"
$ws.Range("C326").Value = "According to Document 1, a curve can wrap a maximum of 50 times."

$ws.Range("A327").Value = "How many curves can be plotted in GEO?"
$ws.Range("B327").Value = "llama3.2:latest"
$ws.Range("C327").Value = "According to Document 1, a curve can wrap a maximum of 50 times."

$ws.Range("A328").Value = "How many curves can be plotted in GEO?"
$ws.Range("B328").Value = "llama3.2:latest"
$ws.Range("C328").Value = "According to Document 1, a curve can wrap a maximum of 50 times."

$ws.Range("A329").Value = "What is the maximum number of data points allowed per curve?"
$ws.Range("B329").Value = "llama3.2:latest"
$ws.Range("C329").Value = "The maximum number of data points allowed per curve is unlimited."

$ws.Range("A330").Value = "What is the maximum number of data points allowed per curve?"
$ws.Range("B330").Value = "llama3.2:latest"
$ws.Range("C330").Value = "The maximum number of data points allowed per curve is unlimited."

$ws.Range("A331").Value = "Why can't I add 251 curve shades to my log?"
$ws.Range("B331").Value = "llama3.2:latest"
$ws.Range("C331").Value = "According to the provided documents, the limit for adding curve shades is not explicitly stated. However, in the section ""Apply lithology shading based on the curve values"", it is mentioned that you should select an unused curve number from CurveShading. This implies that there might be a limit on the number of available curve numbers.`nUnfortunately, without further information or clarification, I couldn't determine the exact limit for adding curve shades."
